$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Small text fixes inside the "Data:" section's final paragraph
#    - "helps to predict" -> "help to predict"
#    - "desc,light"       -> "desc, light"   (comma+space instead of comma)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("helps to predict", $true, $false, $false, $false, $false, $true, 1, $false, "help to predict", 2) | Out-Null
$d.Content.Find.Execute(",light", $true, $false, $false, $false, $false, $true, 1, $false, ", light", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Move the hidden "_GoBack" bookmark from the end of that paragraph to the
#    end of the (new) last real content paragraph ("Conclusion" paragraph),
#    matching where the authors left off editing in the final version.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3) Find the paragraph that now ends with "...build the model." (still the
#    6th paragraph) and append all of the new report sections after it.
# ---------------------------------------------------------------------------
$dataPara = $d.Paragraphs(6)
$dataPara.Range.InsertParagraphAfter()

$newTexts = @(
    "Methodology:",
    "`tThe data has been analyzed with various parameters; like some of the columns have unknown values and it might affect the model fitness. The data has been plotted, against severity code and weather conditions, most of the accidents took place in wet conditions. Like said above, data is being plotted against multiple factors, road conditions, light conditions etc. The chosen test data has multiple factors to train the model and to predict the severity.",
    "`tModel has been built based on classification algorithms. KNN, Decision tree, Support Vector Machine, Logistic Regression algorithms has been used, and the best performance algorithm is decided based on the jaccard score and F1 score. ",
    "Based on the metrics, all algorithms mostly returned with similar score, the data looks biased. The data has to be chosen with multiple variations to find the best model. So the data has been chosen with multiple limits and combinations, the model have been built on. ",
    "Results:",
    "`tThe metrics are being calculated and predicted results are being displayed in the notebook. Based on the metrics , Support vector machine model is best performing having the jaccard score of 1.",
    "Observations:",
    "`tSince the sample dataset is large, it’s better to split and build the model with different parameters.",
    "Conclusion:",
    "`tThe car accident severity has been predicted using classification models."
)

$paraIndex = 7
foreach ($t in $newTexts) {
    $p = $d.Paragraphs($paraIndex)
    $p.Range.Text = $t
    $p.Range.InsertParagraphAfter()
    $paraIndex = $paraIndex + 1
}

# Re-insert the "_GoBack" bookmark at the end of the "Conclusion" paragraph's text
# (paragraph just before the one we are about to keep inserting empties into).
$conclusionPara = $d.Paragraphs($paraIndex - 1)
$d.Bookmarks.Add("_GoBack", $conclusionPara.Range) | Out-Null

# ---------------------------------------------------------------------------
# 4) Two blank paragraphs, then a paragraph containing only a tab.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs($paraIndex)
$p.Range.InsertParagraphAfter()
$paraIndex = $paraIndex + 1

$p = $d.Paragraphs($paraIndex)
$p.Range.Text = "`t"
